# Refresh the market-board derived profit columns (H:N) across the
# item-crafting sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) to match the
# latest scheduled-runner price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1646.4286
$ws.Range("I40").Value = 1312.5
$ws.Range("J40").Value = 1780
$ws.Range("K40").Value = 1312.5
$ws.Range("L40").Value = 1780
$ws.Range("M40").Value = -1137.5
$ws.Range("N40").Value = -2130

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2932.6667
$ws.Range("I64").Value = 2400
$ws.Range("J64").Value = 3287.7778
$ws.Range("K64").Value = 2400
$ws.Range("L64").Value = 3287.7778
$ws.Range("M64").Value = -2152
$ws.Range("N64").Value = -3783.7778

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 2932.6667
$ws.Range("I67").Value = 2400
$ws.Range("J67").Value = 3287.7778
$ws.Range("K67").Value = 2400
$ws.Range("L67").Value = 3287.7778
$ws.Range("M67").Value = -1542
$ws.Range("N67").Value = -5003.7778

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5002.2593
$ws.Range("I74").Value = 7245.5454
$ws.Range("J74").Value = 3460
$ws.Range("K74").Value = 7245.5454
$ws.Range("L74").Value = 3460
$ws.Range("M74").Value = -6309.5454
$ws.Range("N74").Value = -5332

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 5002.2593
$ws.Range("I77").Value = 7245.5454
$ws.Range("J77").Value = 3460
$ws.Range("K77").Value = 36227.727
$ws.Range("L77").Value = 17300
$ws.Range("M77").Value = -31547.727
$ws.Range("N77").Value = -26660

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1830.7084

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1386.1875
$ws.Range("I61").Value = 857.6061
$ws.Range("J61").Value = 2549.0667
$ws.Range("K61").Value = 857.6061
$ws.Range("L61").Value = 2549.0667
$ws.Range("M61").Value = -645.6061
$ws.Range("N61").Value = -2973.0667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2501289.8
$ws.Range("I63").Value = 2501289.8
$ws.Range("K63").Value = 2501289.8
$ws.Range("M63").Value = -2500603.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2501289.8
$ws.Range("I66").Value = 2501289.8
$ws.Range("K66").Value = 12506449
$ws.Range("M66").Value = -12503017

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H113").Value = 59398
$ws.Range("J113").Value = 59398
$ws.Range("L113").Value = 59398
$ws.Range("N113").Value = -68076

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1386.1875
$ws.Range("I136").Value = 857.6061
$ws.Range("J136").Value = 2549.0667
$ws.Range("K136").Value = 2572.8183
$ws.Range("L136").Value = 7647.2001
$ws.Range("M136").Value = -22.81829999999991
$ws.Range("N136").Value = -12747.2001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 811.2222
$ws.Range("I107").Value = 466.83334
$ws.Range("K107").Value = 466.83334
$ws.Range("M107").Value = 1453.16666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 24237.467
$ws.Range("I134").Value = 30314.086
$ws.Range("J134").Value = 2969.3
$ws.Range("K134").Value = 90942.258
$ws.Range("L134").Value = 8907.900000000001
$ws.Range("M134").Value = -88407.258
$ws.Range("N134").Value = -13977.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 393.65216
$ws.Range("I22").Value = 340.3846
$ws.Range("J22").Value = 462.9
$ws.Range("K22").Value = 340.3846
$ws.Range("L22").Value = 462.9
$ws.Range("M22").Value = 9.615400000000022
$ws.Range("N22").Value = -1162.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8334199
$ws.Range("I58").Value = 834.9375
$ws.Range("K58").Value = 834.9375
$ws.Range("M58").Value = -631.9375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2278.5483
$ws.Range("I132").Value = 2252.84
$ws.Range("J132").Value = 2385.6667
$ws.Range("K132").Value = 6758.52
$ws.Range("L132").Value = 7157.000100000001
$ws.Range("M132").Value = -4228.52
$ws.Range("N132").Value = -12217.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1400.1538
$ws.Range("I134").Value = 1336.16
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 4008.48
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -1473.48
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 8334199
$ws.Range("I136").Value = 834.9375
$ws.Range("K136").Value = 2504.8125
$ws.Range("M136").Value = 45.1875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1149.5769
$ws.Range("I5").Value = 345.63635
$ws.Range("J5").Value = 1739.1333
$ws.Range("K5").Value = 1036.90905
$ws.Range("L5").Value = 5217.3999
$ws.Range("M5").Value = -924.90905
$ws.Range("N5").Value = -5441.3999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1027.2963
$ws.Range("I122").Value = 1063.6923
$ws.Range("J122").Value = 993.5
$ws.Range("K122").Value = 9573.2307
$ws.Range("L122").Value = 8941.5
$ws.Range("M122").Value = -7123.2307
$ws.Range("N122").Value = -13841.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2055.077
$ws.Range("I132").Value = 1214.7778
$ws.Range("J132").Value = 2499.9412
$ws.Range("K132").Value = 10933.0002
$ws.Range("L132").Value = 22499.4708
$ws.Range("M132").Value = -8403.0002
$ws.Range("N132").Value = -27559.4708

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1149.5769
$ws.Range("I135").Value = 345.63635
$ws.Range("J135").Value = 1739.1333
$ws.Range("K135").Value = 3110.72715
$ws.Range("L135").Value = 15652.1997
$ws.Range("M135").Value = -575.7271499999997
$ws.Range("N135").Value = -20722.1997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 34778920
$ws.Range("I137").Value = 23811330
$ws.Range("J137").Value = 50133544
$ws.Range("K137").Value = 71433990
$ws.Range("L137").Value = 150400632
$ws.Range("M137").Value = -71428890
$ws.Range("N137").Value = -150410832

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 2276.6667
$ws.Range("I138").Value = 2276.6667
$ws.Range("K138").Value = 6830.000100000001
$ws.Range("M138").Value = -1690.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 3860
$ws.Range("I139").Value = 3633.3333
$ws.Range("J139").Value = 4200
$ws.Range("K139").Value = 10899.9999
$ws.Range("L139").Value = 12600
$ws.Range("M139").Value = -5759.999899999999
$ws.Range("N139").Value = -22880

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2435.9333
$ws.Range("I140").Value = 1730.8182
$ws.Range("J140").Value = 4375
$ws.Range("K140").Value = 5192.4546
$ws.Range("L140").Value = 13125
$ws.Range("M140").Value = -12.45460000000003
$ws.Range("N140").Value = -23485

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 2252.8462
$ws.Range("I141").Value = 2155.3914
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 6466.174199999999
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -1286.174199999999
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 35000
$ws.Range("J114").Value = 35000
$ws.Range("L114").Value = 35000
$ws.Range("N114").Value = -43678

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 33480.125
$ws.Range("I132").Value = 40468.617
$ws.Range("K132").Value = 121405.851
$ws.Range("M132").Value = -118875.851

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6673.921
$ws.Range("I132").Value = 8858.423000000001
$ws.Range("J132").Value = 1940.8334
$ws.Range("K132").Value = 26575.269
$ws.Range("L132").Value = 5822.5002
$ws.Range("M132").Value = -24045.269
$ws.Range("N132").Value = -10882.5002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5180.1113
$ws.Range("I136").Value = 5731.5713
$ws.Range("J136").Value = 3250
$ws.Range("K136").Value = 17194.7139
$ws.Range("L136").Value = 9750
$ws.Range("M136").Value = -14644.7139
$ws.Range("N136").Value = -14850

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1274.4103
$ws.Range("I132").Value = 1239.2122
$ws.Range("J132").Value = 1468
$ws.Range("K132").Value = 3717.6366
$ws.Range("L132").Value = 4404
$ws.Range("M132").Value = -1187.6366
$ws.Range("N132").Value = -9464

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2895.4746
$ws.Range("I136").Value = 3142.6123
$ws.Range("J136").Value = 1684.5
$ws.Range("K136").Value = 9427.836899999998
$ws.Range("L136").Value = 5053.5
$ws.Range("M136").Value = -6877.836899999998
$ws.Range("N136").Value = -10153.5

